# Fixed up regexes and divide by zero crash
#
# The previous run of the simulation divided by zero somewhere in the
# block-size > 1 path, which produced bogus (unnormalized / inflated)
# Throughput and Average-Transmissions-Per-Frame numbers. After fixing the
# regex parsing and the divide-by-zero bug, the simulation was re-run and
# produced corrected values for rows 2-20 (columns A, B, C = Throughput /
# left / right interval, columns E, F, G = Average Transmissions Per Frame /
# left / right interval). Column/row headers in row 1 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.39392000436782837
$ws.Range("B2").Value = 0.37689036130905151
$ws.Range("C2").Value = 0.41094964742660522
$ws.Range("E2").Value = 2.8668475151062012
$ws.Range("F2").Value = 2.6732652187347412
$ws.Range("G2").Value = 3.0604300498962402

$ws.Range("A3").Value = 0.52351999282836914
$ws.Range("B3").Value = 0.50625449419021606
$ws.Range("C3").Value = 0.54078549146652222
$ws.Range("E3").Value = 2.0024442672729492
$ws.Range("F3").Value = 1.9128973484039307
$ws.Range("G3").Value = 2.0919909477233887

$ws.Range("A4").Value = 0.60048002004623413
$ws.Range("B4").Value = 0.58450525999069214
$ws.Range("C4").Value = 0.61645478010177612
$ws.Range("E4").Value = 1.6865731477737427
$ws.Range("F4").Value = 1.6304736137390137
$ws.Range("G4").Value = 1.7426726818084717

$ws.Range("A5").Value = 0.65087997913360596
$ws.Range("B5").Value = 0.63519388437271118
$ws.Range("C5").Value = 0.66656613349914551
$ws.Range("E5").Value = 1.5434224605560303
$ws.Range("F5").Value = 1.4946485757827759
$ws.Range("G5").Value = 1.5921964645385742

$ws.Range("A6").Value = 0.69488000869750977
$ws.Range("B6").Value = 0.67959630489349365
$ws.Range("C6").Value = 0.71016371250152588
$ws.Range("E6").Value = 1.4430468082427979
$ws.Range("F6").Value = 1.391973614692688
$ws.Range("G6").Value = 1.4941198825836182

$ws.Range("A7").Value = 0.72703999280929565
$ws.Range("B7").Value = 0.71297943592071533
$ws.Range("C7").Value = 0.74110054969787598
$ws.Range("E7").Value = 1.361954927444458
$ws.Range("F7").Value = 1.3231201171875
$ws.Range("G7").Value = 1.400789737701416

$ws.Range("A8").Value = 0.74720001220703125
$ws.Range("B8").Value = 0.73356723785400391
$ws.Range("C8").Value = 0.76083278656005859
$ws.Range("E8").Value = 1.3170207738876343
$ws.Range("F8").Value = 1.2886656522750854
$ws.Range("G8").Value = 1.3453758955001831

$ws.Range("A9").Value = 0.77296000719070435
$ws.Range("B9").Value = 0.75990742444992065
$ws.Range("C9").Value = 0.78601258993148804
$ws.Range("E9").Value = 1.2685384750366211
$ws.Range("F9").Value = 1.243529200553894
$ws.Range("G9").Value = 1.2935478687286377

$ws.Range("A10").Value = 0.78655999898910522
$ws.Range("B10").Value = 0.77362906932830811
$ws.Range("C10").Value = 0.79949092864990234
$ws.Range("E10").Value = 1.2452138662338257
$ws.Range("F10").Value = 1.2213971614837646
$ws.Range("G10").Value = 1.2690305709838867

$ws.Range("A11").Value = 0.80111998319625854
$ws.Range("B11").Value = 0.78844368457794189
$ws.Range("C11").Value = 0.81379634141921997
$ws.Range("E11").Value = 1.2233299016952515
$ws.Range("F11").Value = 1.1964647769927979
$ws.Range("G11").Value = 1.2501950263977051

$ws.Range("A12").Value = 0.81455999612808228
$ws.Range("B12").Value = 0.80184060335159302
$ws.Range("C12").Value = 0.82727938890457153
$ws.Range("E12").Value = 1.2007133960723877
$ws.Range("F12").Value = 1.1781561374664307
$ws.Range("G12").Value = 1.2232707738876343

$ws.Range("A13").Value = 0.82543998956680298
$ws.Range("B13").Value = 0.81324422359466553
$ws.Range("C13").Value = 0.83763575553894043
$ws.Range("E13").Value = 1.1821056604385376
$ws.Range("F13").Value = 1.1616504192352295
$ws.Range("G13").Value = 1.2025609016418457

$ws.Range("A14").Value = 0.83407998085021973
$ws.Range("B14").Value = 0.82244986295700073
$ws.Range("C14").Value = 0.8457101583480835
$ws.Range("E14").Value = 1.1679757833480835
$ws.Range("F14").Value = 1.1486247777938843
$ws.Range("G14").Value = 1.1873267889022827

$ws.Range("A15").Value = 0.84112000465393066
$ws.Range("B15").Value = 0.82982927560806274
$ws.Range("C15").Value = 0.85241073369979858
$ws.Range("E15").Value = 1.1564970016479492
$ws.Range("F15").Value = 1.1387484073638916
$ws.Range("G15").Value = 1.1742455959320068

$ws.Range("A16").Value = 0.84975999593734741
$ws.Range("B16").Value = 0.83901971578598022
$ws.Range("C16").Value = 0.8605002760887146
$ws.Range("E16").Value = 1.1430355310440063
$ws.Range("F16").Value = 1.1264986991882324
$ws.Range("G16").Value = 1.1595723628997803

$ws.Range("A17").Value = 0.85343998670578003
$ws.Range("B17").Value = 0.84287291765213013
$ws.Range("C17").Value = 0.86400711536407471
$ws.Range("E17").Value = 1.1372554302215576
$ws.Range("F17").Value = 1.1216062307357788
$ws.Range("G17").Value = 1.1529045104980469

$ws.Range("A18").Value = 0.85887998342514038
$ws.Range("B18").Value = 0.8479427695274353
$ws.Range("C18").Value = 0.86981719732284546
$ws.Range("E18").Value = 1.1309342384338379
$ws.Range("F18").Value = 1.1147351264953613
$ws.Range("G18").Value = 1.1471333503723145

$ws.Range("A19").Value = 0.86176002025604248
$ws.Range("B19").Value = 0.85158246755599976
$ws.Range("C19").Value = 0.87193757295608521
$ws.Range("E19").Value = 1.1252796649932861
$ws.Range("F19").Value = 1.1103427410125732
$ws.Range("G19").Value = 1.140216588973999

$ws.Range("A20").Value = 0.86303997039794922
$ws.Range("B20").Value = 0.85273200273513794
$ws.Range("C20").Value = 0.87334799766540527
$ws.Range("E20").Value = 1.1241039037704468
$ws.Range("F20").Value = 1.1086244583129883
$ws.Range("G20").Value = 1.1395833492279053
